$d = $word.ActiveDocument

# Locate the word "из" inside "...те из можно не учитывать ..." so we can
# fix the typo "из" -> "их" (replacing only the final letter, "з" -> "х"),
# matching how the author actually retyped the single character.
$findRng = $d.Content
[void]$findRng.Find.Execute("те из можно", $false)
$zStart = $findRng.Start + 4
$zEnd = $zStart + 1

# Sanity check: the located single character must be "з".
$checkRng = $d.Range($zStart, $zEnd)
if ($checkRng.Text -ne "з") {
    throw "Could not locate the character to replace"
}

# Insert a temporary marker right before the "з" so the engine keeps the
# text preceding it (".. те и") as its own run once we edit the text that
# follows.
$d.Bookmarks.Add("TempSplitMark", $d.Range($zStart, $zStart))

# Replace "з" with "х" (из -> их).
$editRng = $d.Range($zStart, $zEnd)
$editRng.Text = "х"

# Remove the temporary marker now that the run split it introduced is in
# place.
$d.Bookmarks("TempSplitMark").Delete()

# Re-create the "_GoBack" bookmark (Word's "last edit location" marker)
# right after the newly typed "х" -- this both records the edit point and
# removes the previous "_GoBack" bookmark further up the document, since
# Word keeps only a single "_GoBack" bookmark at a time.
$d.Bookmarks.Add("_GoBack", $d.Range($zStart + 1, $zStart + 1))
